$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-48.
# Some new Price values are plain decimals (e.g. "0.9994") that Excel would
# otherwise auto-convert to a Number; the source data keeps these as plain
# text (same as the "multi-dot" values like "29.525.10" already are), so we
# force text format ("@") on just those cells before assigning the value.
$ws.Range("D2").Value = "29.525.10"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.850.08"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.68"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6318"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.94"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07549"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07676"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.884.35"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6852"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.74"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009839"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "2.133.75"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.214"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "29.575.71"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "234.12"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.608"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.84"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1389"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.72"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.482"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -6.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.283"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.041"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.896"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7163"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.591"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "1.240.48"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01775"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9137"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.139"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "2.041.47"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.97"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.49"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.278"
$ws.Range("E48").Value = "  +9.20%  "

# Rows 49 and 50 swap coin data (EnergySwap <-> BabyDogeCoin) with updated price/volume
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000118"
$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.148"
$ws.Range("E50").Value = "  -0.08%  "
